# The <id> tag in the "p006r_1" paragraph was typed across three separate
# runs (one for "<id>", one for "p006r_1", one for "</id>") each with its
# own run formatting. Collapse them into a single run containing the full
# "<id>p006r_1</id>" text, keeping the Courier New / brown run formatting
# that the "<id>" and "</id>" runs already used.

$d = $word.ActiveDocument

# Locate the opening "<id>" run and the closing "</id>" run by searching
# the document text (robust to any offset/ordering assumptions).
$openTag = $d.Content
$openTag.Find.Execute("<id>", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$closeTag = $d.Content
$closeTag.Find.Execute("</id>", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Everything from right after "<id>" through the end of "</id>" is the
# "p006r_1</id>" text that currently lives in two extra runs; remove it …
$tail = $d.Range($openTag.End, $closeTag.End)
$tail.Delete()

# … and retype it immediately after the surviving "<id>" run so it merges
# into that run and inherits its (Courier New, brown, 18) formatting.
$openTag2 = $d.Range($openTag.Start, $openTag.End)
$openTag2.InsertAfter("p006r_1</id>")
